# Updated main GSC export data.
#
# A new day's row (2025-11-06) arrived at the top of the daily Coverage
# export. The previous export's oldest day (2025-11-05, row 2 on the
# "Chart" sheet) is dropped and every following day's data shifts up by
# one row, so each date now holds the figures that used to belong to the
# next day. The table shrinks from A1:D89 down to A1:D88.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Deleting the entire row 2 shifts rows 3..89 up to 2..88 automatically,
# which is exactly the "drop oldest day, move the rest up" update.
$ws.Rows.Item(2).Delete()
